$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (Beta)
$ws.Range("F2").Value = 286.6905359754097
$ws.Range("G2").Value = 22.58042306437388
$ws.Range("H2").Value = 564.0348173645642
$ws.Range("I2").Value = 1.391645043379825
$ws.Range("J2").Value = 0.4876788257107756
$ws.Range("K2").Value = 2.532756073445229
$ws.Range("L2").Value = 0.2143236358483523
$ws.Range("M2").Value = 0.02909335247701602
$ws.Range("N2").Value = 0.442366767477136

# Row 3 (Gamma)
$ws.Range("F3").Value = 0.02057783878571726
$ws.Range("G3").Value = 0.01263059561467532
$ws.Range("H3").Value = 0.02826543190683595
$ws.Range("I3").Value = 0.01915870218695753
$ws.Range("J3").Value = 0.01165213413616959
$ws.Range("K3").Value = 0.02641194835243531
$ws.Range("L3").Value = 0.020491102423409
$ws.Range("M3").Value = 0.01255505997234859
$ws.Range("N3").Value = 0.02816958245905602

# Row 4 (Beta + Gamma)
$ws.Range("F4").Value = 286.7111138141954
$ws.Range("G4").Value = 22.59305365998856
$ws.Range("H4").Value = 564.063082796471
$ws.Range("I4").Value = 1.410803745566783
$ws.Range("J4").Value = 0.4993309598469452
$ws.Range("K4").Value = 2.559168021797665
$ws.Range("L4").Value = 0.2348147382717614
$ws.Range("M4").Value = 0.04164841244936461
$ws.Range("N4").Value = 0.470536349936192

$wb.Save()
